$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grading data for Essay89-Essay100 (rows 82-93)
$ws.Range("B82").Value = 0
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = 3
$ws.Range("E82").Value = 5
$ws.Range("F82").Value = 0
$ws.Range("H82").Value = "Did not focus on the suggested topic, and without AI statement"

$ws.Range("B83").Value = 2
$ws.Range("C83").Value = 3
$ws.Range("D83").Value = 3
$ws.Range("E83").Value = 5
$ws.Range("F83").Value = 4
$ws.Range("H83").Value = "second part is main generated with AI"

$ws.Range("B84").Value = 2
$ws.Range("C84").Value = 3
$ws.Range("D84").Value = 3
$ws.Range("E84").Value = 3
$ws.Range("F84").Value = 0
$ws.Range("H84").Value = "Focus is not correct; no AI statement; short than required length "

$ws.Range("B85").Value = 3
$ws.Range("C85").Value = 3
$ws.Range("D85").Value = 3
$ws.Range("E85").Value = 3
$ws.Range("F85").Value = 0
$ws.Range("H85").Value = "logic is not good; short; no AI statement"

$ws.Range("B86").Value = 3
$ws.Range("C86").Value = 2
$ws.Range("D86").Value = 3
$ws.Range("E86").Value = 3
$ws.Range("F86").Value = 0
$ws.Range("H86").Value = "No development; no AI statement "

$ws.Range("B87").Value = 2
$ws.Range("C87").Value = 2
$ws.Range("D87").Value = 2
$ws.Range("E87").Value = 2
$ws.Range("F87").Value = 0
$ws.Range("H87").Value = "No development; no AI statement "

$ws.Range("B88").Value = 1
$ws.Range("C88").Value = 2
$ws.Range("D88").Value = 2
$ws.Range("E88").Value = 2
$ws.Range("F88").Value = 0
$ws.Range("H88").Value = "No development; no AI statement "

$ws.Range("B89").Value = 0
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = 2
$ws.Range("E89").Value = 5
$ws.Range("F89").Value = 0
$ws.Range("H89").Value = "Messy logics"

$ws.Range("B90").Value = 2
$ws.Range("C90").Value = 2
$ws.Range("D90").Value = 2
$ws.Range("E90").Value = 4
$ws.Range("F90").Value = 2
$ws.Range("H90").Value = "Did not talk about soluations "

$ws.Range("B91").Value = 2
$ws.Range("C91").Value = 2
$ws.Range("D91").Value = 2
$ws.Range("E91").Value = 3
$ws.Range("F91").Value = 2
$ws.Range("H91").Value = "Vague. Just one paragraph"

$ws.Range("B92").Value = 1
$ws.Range("C92").Value = 2
$ws.Range("D92").Value = 2
$ws.Range("E92").Value = 1
$ws.Range("F92").Value = 0
$ws.Range("H92").Value = "Poor logics"

$ws.Range("B93").Value = 1
$ws.Range("C93").Value = 1
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = 1
$ws.Range("F93").Value = 0
$ws.Range("H93").Value = "too short; has little to do with the required topic"

# Extend the SUM formula down through row 93 (matches the shared-formula group)
$ws.Range("G82:G93").Formula = "=SUM(B82:F82)"

# Rows 91-93 previously had no G cell at all; give them the same centered
# formatting used throughout the rest of the G column (style of G82) without
# disturbing their formulas/values
$ws.Range("G82").Copy()
$ws.Range("G91:G93").PasteSpecial(-4122)

# Match the style used for N1 (same visual style, duplicate xf removed)
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").HorizontalAlignment = -4108
$ws.Range("N1").VerticalAlignment = -4108

# Restore the originally selected cell / scroll position
$ws.Range("I96").Select()
